# refactor: change to PSA naming
# Rename the CGE sector/commodity codes used throughout the "variables" sheets.
# Plain labels (e.g. "AGR") and tuple-style labels (e.g. "('CAP', 'AGR')")
# both get their recognised tokens swapped for the new PSA names; tokens
# that are not in the map (e.g. CAP, LAB, PBS, PAD) are left untouched.

$map = @{
    "AGR" = "AFF"
    "MIN" = "MAQ"
    "MAN" = "MFG"
    "ESW" = "ESWW"
    "CON" = "CNS"
    "WRT" = "TRD"
    "TRS" = "TAS"
    "AFS" = "AFSA"
    "INF" = "IAC"
    "FIN" = "FIA"
    "REA" = "REOD"
    "EDU" = "EDUC"
    "HHS" = "HHSW"
    "OTH" = "OS"
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $old = $cell.Value()

        if ($old -eq $null) {
            continue
        }

        if (-not ($old -is [string])) {
            continue
        }

        if ($old -match "^\('([A-Za-z0-9]+)', '([A-Za-z0-9]+)'\)$") {
            $first = $matches[1]
            $second = $matches[2]

            $newFirst = $first
            if ($map.ContainsKey($first)) {
                $newFirst = $map[$first]
            }

            $newSecond = $second
            if ($map.ContainsKey($second)) {
                $newSecond = $map[$second]
            }

            if (($newFirst -ne $first) -or ($newSecond -ne $second)) {
                $cell.Value = "('$newFirst', '$newSecond')"
            }
        }
        elseif ($map.ContainsKey($old)) {
            $cell.Value = $map[$old]
        }
    }
}
